$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new columns starting at column B, shifting the existing
# B,C,D,E columns (previous date columns + ratings) to E,F,G,H.
$ws.Range("B:D").Insert()

# Keep the fixed-width look of the date columns (C..H), matching the
# width previously used for the date columns before the insert.
$ws.Columns.Item(3).ColumnWidth = 7.1
$ws.Columns.Item(4).ColumnWidth = 7.1
$ws.Columns.Item(5).ColumnWidth = 7.1
$ws.Columns.Item(6).ColumnWidth = 7.1
$ws.Columns.Item(7).ColumnWidth = 7.1
$ws.Columns.Item(8).ColumnWidth = 7.1

# New header row values for the 2 newly added date columns plus the
# freshly opened latest-date column.
$ws.Range("B1").Value = "Jun_27"
$ws.Range("C1").Value = "Jun_26"
$ws.Range("D1").Value = "Jun_26"

# Fill the newly inserted (currently blank) cells for every analyst row
# with the default "UN" (unchanged) marker, matching column B.
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 2).Value = "UN"
    $ws.Cells.Item($r, 3).Value = "UN"
    $ws.Cells.Item($r, 4).Value = "UN"
}

# Add two new analyst rows at the bottom of the table.
$ws.Range("A28").Value = "Benchmark"
$ws.Range("B28").Value = "UN"
$ws.Range("C28").Value = "UN"
$ws.Range("D28").Value = "UN"

$ws.Range("A29").Value = "Evercore ISI"
$ws.Range("B29").Value = "UN"
$ws.Range("C29").Value = "UN"
$ws.Range("D29").Value = "UN"
